$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage (avoids Excel's
# automatic number/date inference silently mangling values such as
# '10.80' -> 10.8 or '58.590.55' style strings).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = 'Normal'
}

Set-TextValue 'D2' '58.284.73'
Set-TextValue 'E2' '  -3.85%  '
Set-TextValue 'D3' '2.742.46'
Set-TextValue 'E3' '  -5.78%  '
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '495.48'
Set-TextValue 'E5' '  -5.92%  '
Set-TextValue 'D6' '133.84'
Set-TextValue 'E6' '  -6.85%  '
Set-TextValue 'E7' '  -0.10%  '
Set-TextValue 'D8' '0.524'
Set-TextValue 'E8' '  -5.34%  '
Set-TextValue 'D9' '2.745.56'
Set-TextValue 'E9' '  -5.58%  '
Set-TextValue 'D10' '5.86'
Set-TextValue 'E10' '  +0.34%  '
Set-TextValue 'E11' '  -5.66%  '
Set-TextValue 'D12' '0.342'
Set-TextValue 'E12' '  -2.79%  '
Set-TextValue 'E13' '  +1.23%  '
Set-TextValue 'D14' '3.233.03'
Set-TextValue 'E14' '  -5.61%  '
Set-TextValue 'D15' '58.481.39'
Set-TextValue 'E15' '  -4.87%  '
Set-TextValue 'D16' '21.13'
Set-TextValue 'E16' '  -6.92%  '
Set-TextValue 'D17' '2.760.75'
Set-TextValue 'E17' '  -5.34%  '
Set-TextValue 'D18' '0.0000132'
Set-TextValue 'E18' '  -5.44%  '
Set-TextValue 'D19' '4.63'
Set-TextValue 'E19' '  -6.23%  '
Set-TextValue 'D20' '344.23'
Set-TextValue 'E20' '  -4.30%  '
Set-TextValue 'D21' '10.80'
Set-TextValue 'E21' '  -6.31%  '
Set-TextValue 'D22' '6.11'
Set-TextValue 'E22' '  -5.20%  '
Set-TextValue 'D23' '0.997'
Set-TextValue 'E23' '  -0.24%  '
Set-TextValue 'D24' '5.61'
Set-TextValue 'E24' '  -1.03%  '
Set-TextValue 'D25' '62.07'
Set-TextValue 'E25' '  -2.46%  '
Set-TextValue 'D26' '0.419'
Set-TextValue 'E26' '  -6.76%  '
Set-TextValue 'D27' '0.170'
Set-TextValue 'E27' '  -7.34%  '
Set-TextValue 'E28' '  +0.12%  '
Set-TextValue 'D29' '7.22'
Set-TextValue 'E29' '  -5.53%  '
Set-TextValue 'D30' '0.0₃0785'
Set-TextValue 'E30' '  -8.29%  '
Set-TextValue 'D31' '0.999'
Set-TextValue 'E31' '  -0.05%  '
Set-TextValue 'D32' '1.58'
Set-TextValue 'D33' '18.82'
Set-TextValue 'E33' '  -4.02%  '
Set-TextValue 'D34' '148.72'
Set-TextValue 'E34' '  -3.24%  '
Set-TextValue 'D35' '4.09'
Set-TextValue 'E35' '  -5.72%  '
Set-TextValue 'D36' '5.24'
Set-TextValue 'E36' '  -6.14%  '
Set-TextValue 'D37' '0.892'
Set-TextValue 'E37' '  -10.95%  '
Set-TextValue 'B38' 'ImmutableX'
Set-TextValue 'C38' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D38' '1.11'
Set-TextValue 'E38' '  -7.83%  '
Set-TextValue 'B39' 'OKB'
Set-TextValue 'C39' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D39' '36.58'
Set-TextValue 'E39' '  -3.67%  '
Set-TextValue 'D40' '2.179.07'
Set-TextValue 'E40' '  -6.82%  '
Set-TextValue 'B41' 'FirstDigitalUSD'
Set-TextValue 'C41' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D41' '0.999'
Set-TextValue 'E41' '  -0.09%  '
Set-TextValue 'B42' 'Filecoin'
Set-TextValue 'C42' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D42' '3.47'
Set-TextValue 'E42' '  -5.89%  '
Set-TextValue 'B43' 'Mantle'
Set-TextValue 'C43' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D43' '0.604'
Set-TextValue 'E43' '  -6.02%  '
Set-TextValue 'D44' '0.0548'
Set-TextValue 'E44' '  -3.36%  '
Set-TextValue 'D45' '1.33'
Set-TextValue 'E45' '  -9.70%  '
Set-TextValue 'D46' '18.80'
Set-TextValue 'E46' '  -9.31%  '
Set-TextValue 'D47' '10.34'
Set-TextValue 'E47' '  -0.09%  '
Set-TextValue 'E48' '  -4.41%  '
Set-TextValue 'B49' 'RenderToken'
Set-TextValue 'C49' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D49' '4.51'
Set-TextValue 'E49' '  -6.89%  '
Set-TextValue 'B50' 'Stellar'
Set-TextValue 'C50' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D50' '0.0879'
Set-TextValue 'E50' '  -4.71%  '
Set-TextValue 'D51' '17.17'
Set-TextValue 'E51' '  -6.18%  '
